$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.473.78'
$ws.Range("E2").Value = '  +1.79%  '

$ws.Range("D3").Value = '1.859.25'
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4766'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3787'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07317'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9295'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.68'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07787'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.64%  '

$ws.Range("D13").Value = '1.862.77'
$ws.Range("E13").Value = '  +0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.448'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.562'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.17'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.40%  '

$ws.Range("E17").Value = '  -0.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008824'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.83%  '

$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").Value = '27.393.98'
$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.64'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.090'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("E23").Value = '  +0.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.942'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.12'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.46'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.004'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.31'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.949'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08871'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.330'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7520'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.579'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.703'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02046'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.122'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("E38").Value = '  +5.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05285'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.014'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.556'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1518'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4871'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.67'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.011'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.85'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.664'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06095'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9108'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.00%  '
